$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column F width (matches <col width="26" customWidth="1" min="6" max="6"/>)
# Note: Excel stores column width internally based on pixel-rounding of the
# "characters" width, so 25.14 round-trips to an on-disk width of exactly 26.
$ws.Columns.Item(6).ColumnWidth = 25.14

# Header cell F1 - same style as the other header cells (bold, centered, bordered)
$ws.Range("F1").Value = "Company Verification"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Column F data values, rows 2-31 (matches "Company Verification" column added to diff)
$values = @(
    "FloWorks",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Omnicell",
    "DNOW",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Company does not match",
    "Elbit Systems of America",
    "Company does not match"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
